$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy B1's style (bold, border, centered) into C1 so it matches the header formatting
$ws.Range("B1").Copy($ws.Range("C1"))

# Fill column C (C1:C67) with the new 't+3' predicted values
$values = @(
    2,
    -0.3040417462236097,
    0.2251535355371952,
    0.1544151644258079,
    -0.06922912467150039,
    0.1592649206024101,
    -0.4466575725286778,
    -0.2530181817640409,
    -0.3885918744649107,
    0.3690991585549788,
    -0.2374145889386449,
    -0.06832904299805047,
    -0.01024573843619881,
    0.09767730266494998,
    0.03554015468291617,
    0.4182123011612361,
    0.6463634593515115,
    0.050867864187726,
    0.4182019070579845,
    0.3586001354007108,
    0.5605837047146707,
    0.4216504564096156,
    -0.02984816026661184,
    4.319056669188416,
    0.4598891472662188,
    0.3670980678865449,
    0.2621466084748486,
    0.991487030393535,
    5.224988178542493,
    0.923784986971142,
    -0.2707495582110871,
    0.7504266855551888,
    0.8844208549016394,
    -0.7032468508123622,
    0.800343452240891,
    0.7483092607616741,
    0.7198224738286401,
    0.7053106754767264,
    0.5701106224081278,
    0.7401332273910494,
    0.5508581871494762,
    0.6645571692077971,
    0.6986510744487012,
    0.6638463159649135,
    0.6316393742209037,
    -1.281318576789337,
    -0.9953210513289173,
    -0.8891017246582267,
    -0.6541934193054134,
    -0.05623713313460828,
    -0.88181270688059,
    -0.88181270688059,
    -1.116975934909187,
    -0.1947372290805855,
    -1.015048157973971,
    -0.8977876783020756,
    -0.9300980938665987,
    -1.123637583142322,
    -0.8511685513275964,
    -0.4902850983854158,
    0.3716943380942198,
    -1.203251063982498,
    -0.7335006093035624,
    -0.8991233070265447,
    -0.1057560117175542,
    -0.7847113697185064,
    -0.7583261251229748
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $values[$i]
}
